$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.2933238827060871;  "C2" = 0.2933238827060871
    "B3" = 0.586170229238267;   "C3" = 0.586170229238267
    "B4" = 0.7286494602544143;  "C4" = 0.7286494602544143
    "B5" = 0.0433905709794368;  "C5" = 0.0433905709794368;  "D5" = 191;  "E5" = 191
    "B6" = 0.9432547317543618;  "C6" = 0.9432547317543618
    "B7" = 0.0005065452317266613; "C7" = 0.0005065452317266613; "D7" = 786
    "B8" = 0.1721303739188967;  "C8" = 0.1721303739188967;  "D8" = 142; "E8" = 142
    "B9" = 0.07243033620187289; "C9" = 0.07243033620187289; "D9" = 282; "E9" = 282
    "B10" = 0.4341558216192389; "C10" = 0.4341558216192389; "D10" = 36; "E10" = 36; "F10" = 124
    "B11" = 0.0000006129115011681582; "C11" = 0.0000006129115011681582; "D11" = 1088; "E11" = 1088; "F11" = 1483; "G11" = 1483
    "B12" = 0.1104444481152851; "C12" = 0.1104444481152851
    "B13" = 0.5382430939996513; "C13" = 0.5382430939996513; "D13" = 86; "E13" = 86
    "B14" = 0.7833230176515672; "C14" = 0.7833230176515672; "D14" = 30; "E14" = 30
    "B15" = 0.8094692513399709; "C15" = 0.8094692513399709; "E15" = 18; "F15" = 31
    "B16" = 0.6511559132080488; "C16" = 0.6511559132080488
    "B17" = 0.6585833303580919; "C17" = 0.6585833303580919
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
